# "Prepare for 2022 publications."
# Bump the copyright year on the cover page from 2021 to 2022.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2021", $true, $true, $false, $false, $false, $true, 1, $false, "2022", 2)
